$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.268.47"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.509.74"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'604.50"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'175.37"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("D8").Value = "3.504.06"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("E11").Value = "  +8.00%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'46.23"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "4.067.83"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'611.03"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "3.507.38"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "70.343.92"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "'17.34"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").Value = "'0.876"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -11.21%  "
$ws.Range("D24").Value = "'98.43"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "'15.54"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "'3.72"
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'2.56"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").Value = "'33.84"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").Value = "'9.01"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("B31").Value = "Stacks"
$ws.Range("C31").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D31").Value = "'2.98"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Value = "'1.28"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").Value = "'6.83"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "'623.63"
$ws.Range("E35").Value = "  +12.49%  "
$ws.Range("D36").Value = "'0.0993"
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").Value = "'3.55"
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("D38").Value = "'10.75"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").Value = "'0.0478"
$ws.Range("E39").Value = "  +6.81%  "
$ws.Range("D40").Value = "'56.78"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").Value = "3.369.52"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("D44").Value = "0.0₃0737"
$ws.Range("E44").Value = "  +5.77%  "
$ws.Range("D45").Value = "'0.308"
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").Value = "'32.17"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'2.56"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'0.130"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("D50").Value = "'132.58"
$ws.Range("E50").Value = "  -2.21%  "
